$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "mokhtar"
$ws.Range("F2").Value = 1700203
$ws.Range("F3").Value = 1700204
$ws.Range("F4").Value = 1700205

[void]$ws.Range("G10").Select()
